# Apply the "simulator full-month coverage, persist logs, fix employees" edit.
$wb = $excel.ActiveWorkbook

$tsSheet = $wb.Worksheets.Item("Weekly Timesheet")
$jsSheet = $wb.Worksheets.Item("Jason Schema")

# --- Fix employee / client names (shared strings are reused across sheets) ---
$tsSheet.Range("B2").Value = "Winn"
$tsSheet.Range("B3").Value = "Keevil"
$tsSheet.Range("B4").Value = "Howard"
$tsSheet.Range("B5").Value = "Markfield"
$tsSheet.Range("B6").Value = "Layne"

$jsSheet.Range("D2").Value = "Winn"
$jsSheet.Range("D3").Value = "Keevil"
$jsSheet.Range("D4").Value = "Howard"
$jsSheet.Range("D5").Value = "Markfield"
$jsSheet.Range("D6").Value = "Layne"

# --- Fix employee ID (same ID repeated for every log row) ---
$jsSheet.Range("B2").Value = "emp_75yd72zj"
$jsSheet.Range("B3").Value = "emp_75yd72zj"
$jsSheet.Range("B4").Value = "emp_75yd72zj"
$jsSheet.Range("B5").Value = "emp_75yd72zj"
$jsSheet.Range("B6").Value = "emp_75yd72zj"

# --- Simulator full-month coverage: populate rate/total for each day row ---
$tsSheet.Range("E2").Value = 90
$tsSheet.Range("F2").Value = 720
$tsSheet.Range("E3").Value = 90
$tsSheet.Range("F3").Value = 720
$tsSheet.Range("E4").Value = 90
$tsSheet.Range("F4").Value = 720
$tsSheet.Range("E5").Value = 90
$tsSheet.Range("F5").Value = 720
$tsSheet.Range("E6").Value = 90
$tsSheet.Range("F6").Value = 720

$tsSheet.Range("F8").Value = 3600
$tsSheet.Range("F11").Value = 3600
$tsSheet.Range("F13").Value = 3600

# --- Persist logs: mirror the same rate/total into the Jason Schema log sheet ---
$jsSheet.Range("F2").Value = 90
$jsSheet.Range("G2").Value = 720
$jsSheet.Range("F3").Value = 90
$jsSheet.Range("G3").Value = 720
$jsSheet.Range("F4").Value = 90
$jsSheet.Range("G4").Value = 720
$jsSheet.Range("F5").Value = 90
$jsSheet.Range("G5").Value = 720
$jsSheet.Range("F6").Value = 90
$jsSheet.Range("G6").Value = 720
